$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cfs = $ws.Range("H5,H8").FormatConditions
Write-Host "Count: $($cfs.Count)"
